$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.743.21"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "1.631.92"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.501"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.01"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.53"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0782"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "1.657.15"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "1.857.11"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.552"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.01%  "
$ws.Range("D16").Value = "0.0₃0768"
$ws.Range("E16").Value = "  -0.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.21%  "
$ws.Range("D18").Value = "25.763.21"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -3.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0485"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.64%  "
$ws.Range("E32").Value = "  +0.88%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("E35").Value = "  +0.49%  "
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.547"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").Value = "1.105.15"
$ws.Range("E39").Value = "  -2.12%  "
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.84"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.26%  "
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("D45").Value = "0.0₆0111"
$ws.Range("E45").Value = "  -2.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.02"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("E47").Value = "  -2.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.46%  "
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  +3.01%  "
$ws.Range("E51").Value = "  +0.35%  "
